$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 18:58"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 6351246
$ws.Range("C4").Value = 16002
$ws.Range("D4").Value = 3583205
$ws.Range("E4").Value = 2576649
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 334
$ws.Range("H4").Value = 191392

# Row 5: Brasil -> Brasil
$ws.Range("B5").Value = 4054474
$ws.Range("C5").Value = 8324
$ws.Range("D5").Value = 3247610
$ws.Range("E5").Value = 681942
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 193
$ws.Range("H5").Value = 124922

# Row 6: India -> India
$ws.Range("B6").Value = 4014744
$ws.Range("C6").Value = 81620
$ws.Range("D6").Value = 3101245
$ws.Range("E6").Value = 843885
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1045
$ws.Range("H6").Value = 69614

# Row 12: España -> España
$ws.Range("B12").Value = 498989
$ws.Range("C12").Value = 10476
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 184
$ws.Range("H12").Value = 29418

# Row 14: Chile -> Chile
$ws.Range("B14").Value = 418469
$ws.Range("C14").Value = 1968
$ws.Range("D14").Value = 391248
$ws.Range("E14").Value = 15727
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 72
$ws.Range("H14").Value = 11494

# Row 16: Reino Unido -> Reino Unido
$ws.Range("B16").Value = 342351
$ws.Range("C16").Value = 1940
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 41537

# Row 21: Turquia -> Turquia
$ws.Range("B21").Value = 276555
$ws.Range("C21").Value = 1612
$ws.Range("D21").Value = 249108
$ws.Range("E21").Value = 20883
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = 6564

# Row 22: Italia -> Italia
$ws.Range("B22").Value = 274644
$ws.Range("C22").Value = 1733
$ws.Range("D22").Value = 209027
$ws.Range("E22").Value = 30099
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 35518

# Row 24: Alemania -> Alemania
$ws.Range("B24").Value = 249372
$ws.Range("C24").Value = 558
$ws.Range("D24").Value = 224600
$ws.Range("E24").Value = 15373
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 9399

# Row 28: Canada -> Canada
$ws.Range("B28").Value = 130825
$ws.Range("C28").Value = 332
$ws.Range("D28").Value = 115669
$ws.Range("E28").Value = 6015
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 9141

# Row 29: Israel -> Israel
$ws.Range("B29").Value = 126419
$ws.Range("C29").Value = 1964
$ws.Range("D29").Value = 100357
$ws.Range("E29").Value = 25069
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 993

# Row 55: Barein -> Barein
$ws.Range("B55").Value = 53433
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 50013
$ws.Range("E55").Value = 3225
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 195

# Row 57: Argelia -> Argelia
$ws.Range("B57").Value = 45773
$ws.Range("C57").Value = 304
$ws.Range("D57").Value = 32259
$ws.Range("E57").Value = 11975
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 1539

# Row 70: Irlanda -> Irlanda
$ws.Range("B70").Value = 29303
$ws.Range("C70").Value = 97
$ws.Range("D70").Value = 23364
$ws.Range("E70").Value = 4162
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 1777

# Row 72: Chequia -> Chequia
$ws.Range("B72").Value = 26943
$ws.Range("C72").Value = 491
$ws.Range("D72").Value = 18985
$ws.Range("E72").Value = 7529
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 429

# Row 76: Corea del Sur -> Bosnia y Herzegovina
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 21142
$ws.Range("C76").Value = 338
$ws.Range("D76").Value = 14476
$ws.Range("E76").Value = 6027
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 639

# Row 77: Bosnia y Herzegovina -> Corea del Sur
$ws.Range("A77").Value = "Corea del Sur"
$ws.Range("B77").Value = 20842
$ws.Range("C77").Value = 198
$ws.Range("D77").Value = 15783
$ws.Range("E77").Value = 4728
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 331

# Row 80: Libano -> Libano
$ws.Range("B80").Value = 19490
$ws.Range("C80").Value = 527
$ws.Range("D80").Value = 5592
$ws.Range("E80").Value = 13715
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 4
$ws.Range("H80").Value = 183

# Row 97: Guayana Francesa -> Guayana Francesa
$ws.Range("B97").Value = 9276
$ws.Range("C97").Value = 25
$ws.Range("D97").Value = 8792
$ws.Range("E97").Value = 422
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 62

# Row 139: Aruba -> Jordania
$ws.Range("A139").Value = "Jordania"
$ws.Range("B139").Value = 2301
$ws.Range("C139").Value = 68
$ws.Range("D139").Value = 1676
$ws.Range("E139").Value = 609
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 16

# Row 140: Jordania -> Aruba
$ws.Range("A140").Value = "Aruba"
$ws.Range("B140").Value = 2292
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 1031
$ws.Range("E140").Value = 1248
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 13

# Row 145: Trinidad yTobago -> Trinidad yTobago
$ws.Range("B145").Value = 2030
$ws.Range("C145").Value = 46
$ws.Range("D145").Value = 707
$ws.Range("E145").Value = 1293
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 30

# Row 146: Malta -> Reunion
$ws.Range("A146").Value = "Reunion"
$ws.Range("B146").Value = 2002
$ws.Range("C146").Value = 90
$ws.Range("D146").Value = 880
$ws.Range("E146").Value = 1112
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 10

# Row 147: Yemen -> Malta
$ws.Range("A147").Value = "Malta"
$ws.Range("B147").Value = 1984
$ws.Range("C147").Value = 19
$ws.Range("D147").Value = 1565
$ws.Range("E147").Value = 406
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 13

# Row 148: Reunion -> Yemen
$ws.Range("A148").Value = "Yemen"
$ws.Range("B148").Value = 1979
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 1180
$ws.Range("E148").Value = 228
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 571

# Row 160: Principado de Andorra -> Principado de Andorra
$ws.Range("B160").Value = 1215
$ws.Range("C160").Value = 16
$ws.Range("D160").Value = 928
$ws.Range("E160").Value = 234
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 53
